# Atualização de bases das ligas, do dia: 15-06-2024 às 21:10
#
# For a handful of fixture-row pairs, the data in columns B:AD (match id,
# teams, scores, odds, ...) was swapped between the two rows while the
# leading index column (A) stayed put. Swap each pair back/forth to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC","AD")

$rowPairs = @(
    @(38, 39),
    @(129, 131),
    @(224, 225),
    @(231, 232),
    @(256, 257),
    @(267, 269)
)

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    foreach ($col in $cols) {
        $rng1 = $ws.Range("$col$r1")
        $rng2 = $ws.Range("$col$r2")

        $v1 = $rng1.Value2
        $v2 = $rng2.Value2

        $rng1.Value2 = $v2
        $rng2.Value2 = $v1
    }
}
